# Fruta / hortaliza, semanal
# Insert this week's new price records for "Pera" (Packham's Triumph,
# Especial and Primera) at the top of the Macroferia Regional de Talca
# block, pushing the existing rows down by two positions.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows before row 1110 (formatting of the surrounding
# rows - e.g. the date style on column D - is inherited automatically).
$ws.Rows.Item(1110).Resize(2).Insert()

# New row 1110: Packham's Triumph / Especial
$ws.Cells.Item(1110, 1).Value = 5
$ws.Cells.Item(1110, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(1110, 3).Value = "Maule"
$ws.Cells.Item(1110, 4).Value = 45223
$ws.Cells.Item(1110, 5).Value = 7
$ws.Cells.Item(1110, 6).Value = "Fruta"
$ws.Cells.Item(1110, 7).Value = 100104
$ws.Cells.Item(1110, 8).Value = "Frutos de pepita"
$ws.Cells.Item(1110, 9).Value = 100104005
$ws.Cells.Item(1110, 10).Value = "Pera"
$ws.Cells.Item(1110, 11).Value = "Packham's Triumph"
$ws.Cells.Item(1110, 12).Value = "Especial"
$ws.Cells.Item(1110, 13).Value = 220
$ws.Cells.Item(1110, 14).Value = 16000
$ws.Cells.Item(1110, 15).Value = 16000
$ws.Cells.Item(1110, 16).Value = 16000
$ws.Cells.Item(1110, 17).Value = "`$/bandeja 18 kilos granel"
$ws.Cells.Item(1110, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(1110, 19).Value = 889
$ws.Cells.Item(1110, 20).Value = 18

# New row 1111: Packham's Triumph / Primera
$ws.Cells.Item(1111, 1).Value = 5
$ws.Cells.Item(1111, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(1111, 3).Value = "Maule"
$ws.Cells.Item(1111, 4).Value = 45223
$ws.Cells.Item(1111, 5).Value = 7
$ws.Cells.Item(1111, 6).Value = "Fruta"
$ws.Cells.Item(1111, 7).Value = 100104
$ws.Cells.Item(1111, 8).Value = "Frutos de pepita"
$ws.Cells.Item(1111, 9).Value = 100104005
$ws.Cells.Item(1111, 10).Value = "Pera"
$ws.Cells.Item(1111, 11).Value = "Packham's Triumph"
$ws.Cells.Item(1111, 12).Value = "Primera"
$ws.Cells.Item(1111, 13).Value = 250
$ws.Cells.Item(1111, 14).Value = 14000
$ws.Cells.Item(1111, 15).Value = 14000
$ws.Cells.Item(1111, 16).Value = 14000
$ws.Cells.Item(1111, 17).Value = "`$/bandeja 18 kilos granel"
$ws.Cells.Item(1111, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(1111, 19).Value = 778
$ws.Cells.Item(1111, 20).Value = 18
